$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set E2 and E5 to "NA" (previously empty)
$ws.Range("E2").Value = "NA"
$ws.Range("E5").Value = "NA"

# D10 changes from "in process" to "x"
$ws.Range("D10").Value = "x"

# Update the selection to reflect E11 as the active cell (as in the saved file)
$ws.Range("E11").Select()
